# For every row in the "Recorded By" column (column G) of the active sheet,
# reverse the order of the comma-separated list of recorder names/emails.
# Cells that hold only a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $joined = $reversed -join ", "
        $cell.Value = $joined
    }
}
